# [LC-850] Release of LetsCo OS v1.3.0
# Update the legacy short KPI ids (GP1/GP2/GP3/BP1/BP2/BP3) in column B to
# their zero-padded equivalents (GP01/GP02/GP03/BP01/BP02/BP03), and reset
# the sheet view/selection back to A1 (it used to be parked at G1/I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- id column (B) updates -------------------------------------------------
# Each value below is shared by every row of its KPI group, so every one of
# those rows needs to be rewritten (Excel stores repeated strings once, but
# we drive this purely through the cell Values, same as a user retyping
# them).
$ws.Range("B16").Value = "GP01"        # was GP1
$ws.Range("B17:B18").Value = "GP02"    # was GP2
$ws.Range("B19:B20").Value = "GP03"    # was GP3
$ws.Range("B21:B25").Value = "BP01"    # was BP1
$ws.Range("B26:B30").Value = "BP02"    # was BP2
$ws.Range("B31:B60").Value = "BP03"    # was BP3

# --- reset the view back to A1 ---------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
